# Update the "Input [Source Name]" / "Output [Source Name]" headers to
# "Input [Sample Name]" / "Output [Sample Name]" on the Events-Harvest sheet,
# including the underlying table column definitions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events-Harvest")

# Update the header cell values (these also back the shared strings used by
# the table's column headers, and the table definition's column names are
# kept in sync with these header cells by Excel).
$ws.Range("A1").Value = "Input [Sample Name]"
$ws.Range("AF1").Value = "Output [Sample Name]"
